# Automatische test-sync: 2025-06-19 19:07:30
# Add a new log entry row to the "Logs" sheet and update the
# "Dashboard" summary sheet accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 38 ------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(38, 1).Value = "Sollicitatie marketingfunctie"
$logs.Cells.Item(38, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(38, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Cells.Item(38, 4).Value = "Sollicitatie"
$logs.Cells.Item(38, 6).Value = "2025-06-19 19:07:25"
$logs.Cells.Item(38, 7).Value = "Nee"

# Extend the conditional formatting ranges (D2:D37 -> D2:D38, G2:G37 -> G2:G38)
# so the newly added row participates in the existing category/answer highlighting.
$dConditions = $logs.Range("D2:D37").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D38"))
}

$gConditions = $logs.Range("G2:G37").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G38"))
}

# --- Sheet "Dashboard": update category counts ---------------------------
# The new entry's category ("Sollicitatie") count increases from 3 to 4,
# which also changes the sort order relative to "Bestelling" (3): the two
# rows swap places.
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(5, 1).Value = "Sollicitatie"
$dashboard.Cells.Item(5, 2).Value = 4

$dashboard.Cells.Item(6, 1).Value = "Bestelling"
$dashboard.Cells.Item(6, 2).Value = 3
